$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 94, shifting existing rows 94:206 down to 95:206 (... to 206)
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with the new record
$ws.Cells.Item(94, 1).Value = 7
$ws.Cells.Item(94, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(94, 3).Value = "Ñuble"
$ws.Cells.Item(94, 4).Value = 44546
$ws.Cells.Item(94, 5).Value = 16
$ws.Cells.Item(94, 6).Value = 100112023
$ws.Cells.Item(94, 7).Value = "Brócoli"
$ws.Cells.Item(94, 8).Value = "Sin especificar"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 800
$ws.Cells.Item(94, 11).Value = 700
$ws.Cells.Item(94, 12).Value = 800
$ws.Cells.Item(94, 13).Value = 750
$ws.Cells.Item(94, 14).Value = "$/unidad"
$ws.Cells.Item(94, 15).Value = "Región del Maule"
$ws.Cells.Item(94, 16).Value = 750
$ws.Cells.Item(94, 17).Value = 1
$ws.Cells.Item(94, 18).Value = "Hortaliza"
